# Add instructor/email/position/phone info for newly hired staff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("F4").Style = "Normal"
$ws.Range("F4").Value = 'Chen,Yan'
$ws.Range("G4").Style = "Normal"
$ws.Range("G4").Value = 'ychen@bmcc.cuny.edu'
$ws.Range("H4").Value = 'Faculty'
$ws.Range("I4").Style = "Normal"
$ws.Range("I4").Value = '212-220-8384'

# Row 5
$ws.Range("F5").Style = "Normal"
$ws.Range("F5").Value = 'Melie,Ora K'
$ws.Range("H5").Value = 'Adjunct'
$ws.Range("I5").Style = "Normal"
$ws.Range("I5").Value = '212-220-1476'

# Row 7
$ws.Range("F7").Style = "Normal"
$ws.Range("F7").Value = 'Toliyat Abolhassani,Amir Mohsen'
$ws.Range("H7").Value = 'Adjunct'
$ws.Range("I7").Style = "Normal"
$ws.Range("I7").Value = '212-220-1476'

# Row 8
$ws.Range("F8").Style = "Normal"
$ws.Range("F8").Value = 'Azhar,Mohammad Quamrul'
$ws.Range("G8").Style = "Normal"
$ws.Range("G8").Value = 'mazhar@bmcc.cuny.edu'
$ws.Range("H8").Value = 'Faculty'
$ws.Range("I8").Style = "Normal"
$ws.Range("I8").Value = '212-220-1477'

# Row 9
$ws.Range("F9").Style = "Normal"
$ws.Range("F9").Value = 'Jayaweera,Darshani P'
$ws.Range("H9").Value = 'Adjunct'
$ws.Range("I9").Style = "Normal"
$ws.Range("I9").Value = '212-220-1476'

# Row 11
$ws.Range("F11").Style = "Normal"
$ws.Range("F11").Value = 'Harricharan,Andy L'
$ws.Range("H11").Value = 'Adjunct'
$ws.Range("I11").Style = "Normal"
$ws.Range("I11").Value = '212-220-1476'

# Row 12
$ws.Range("F12").Style = "Normal"
$ws.Range("F12").Value = 'Zeidan,Ayman I'
$ws.Range("H12").Value = 'Adjunct'
$ws.Range("I12").Style = "Normal"
$ws.Range("I12").Value = '212-220-1476'

# Row 13
$ws.Range("F13").Style = "Normal"
$ws.Range("F13").Value = 'Kalicharan,Dharamraj'
$ws.Range("H13").Value = 'Adjunct'
$ws.Range("I13").Style = "Normal"
$ws.Range("I13").Value = '212-220-1476'

# Row 15
$ws.Range("F15").Style = "Normal"
$ws.Range("F15").Value = 'Kalicharan,Dharamraj'
$ws.Range("H15").Value = 'Adjunct'
$ws.Range("I15").Value = '212-220-1476'

# Row 16
$ws.Range("F16").Style = "Normal"
$ws.Range("F16").Value = 'Veyler,Mikhail'
$ws.Range("H16").Value = 'Adjunct'
$ws.Range("I16").Value = '212-220-1476'

# Row 17
$ws.Range("F17").Style = "Normal"
$ws.Range("F17").Value = 'Adesman,Alexander'
$ws.Range("H17").Value = 'Adjunct'
$ws.Range("I17").Value = '212-220-1476'

# Row 18
$ws.Range("F18").Style = "Normal"
$ws.Range("F18").Value = 'Rani,Chigurupati S'
$ws.Range("G18").Style = "Normal"
$ws.Range("G18").Value = 'csrani@bmcc.cuny.edu'
$ws.Range("H18").Value = 'Faculty'
$ws.Range("I18").Style = "Normal"
$ws.Range("I18").Value = '212-220-1487'

# Row 19
$ws.Range("F19").Style = "Normal"
$ws.Range("F19").Value = 'Khan,Lawrence'
$ws.Range("H19").Value = 'Adjunct'
$ws.Range("I19").Value = '212-220-1476'

# Row 21
$ws.Range("F21").Style = "Normal"
$ws.Range("F21").Value = 'Scheiman,Arnold J'
$ws.Range("G21").Style = "Normal"
$ws.Range("G21").Value = 'ascheiman@bmcc.cuny.edu'
$ws.Range("H21").Value = 'Faculty'
$ws.Range("I21").Style = "Normal"
$ws.Range("I21").Value = '212-220-7227'

# Row 22
$ws.Range("F22").Style = "Normal"
$ws.Range("F22").Value = 'Liu,Ligon Mengxu'
$ws.Range("H22").Value = 'Adjunct'
$ws.Range("I22").Value = '212-220-1476'

# Row 23
$ws.Range("F23").Style = "Normal"
$ws.Range("F23").Value = 'Mokal,Prajakta L'
$ws.Range("H23").Value = 'Adjunct'
$ws.Range("I23").Value = '212-220-1476'

# Row 24
$ws.Range("F24").Style = "Normal"
$ws.Range("F24").Value = 'Genis,Yakov'
$ws.Range("G24").Style = "Normal"
$ws.Range("G24").Value = 'ygenis@bmcc.cuny.edu'
$ws.Range("H24").Value = 'Faculty'
$ws.Range("I24").Style = "Normal"
$ws.Range("I24").Value = '212-220-1482'

# Row 25
$ws.Range("F25").Style = "Normal"
$ws.Range("F25").Value = 'Hasan,Naushad'
$ws.Range("H25").Value = 'Adjunct'
$ws.Range("I25").Value = '212-220-1476'

# Row 27
$ws.Range("F27").Style = "Normal"
$ws.Range("F27").Value = 'Vargas,Jose I'
$ws.Range("G27").Style = "Normal"
$ws.Range("G27").Value = 'jvargas@bmcc.cuny.edu'
$ws.Range("H27").Value = 'Faculty'
$ws.Range("I27").Style = "Normal"
$ws.Range("I27").Value = '212-220-1481'

# Row 29
$ws.Range("F29").Style = "Normal"
$ws.Range("F29").Value = 'Cooper,Kenneth J'
$ws.Range("H29").Value = 'Adjunct'
$ws.Range("I29").Value = '212-220-1476'

# Row 33
$ws.Range("F33").Style = "Normal"
$ws.Range("F33").Value = 'Nossa,George A'
$ws.Range("H33").Value = 'Adjunct'
$ws.Range("I33").Value = '212-220-1476'

# Row 34
$ws.Range("F34").Style = "Normal"
$ws.Range("F34").Value = 'O''Faire,Lashawna R'
$ws.Range("H34").Value = 'Adjunct'
$ws.Range("I34").Value = '212-220-1476'

# Row 35
$ws.Range("F35").Style = "Normal"
$ws.Range("F35").Value = 'Lau,Roy'
$ws.Range("H35").Value = 'Adjunct'
$ws.Range("I35").Value = '212-220-1476'

# Row 36
$ws.Range("F36").Style = "Normal"
$ws.Range("F36").Value = 'Hernandez,Candido'
$ws.Range("H36").Value = 'Adjunct'
$ws.Range("I36").Value = '212-220-1476'

# Row 37
$ws.Range("H37").Value = 'Adjunct'
$ws.Range("I37").Value = '212-220-1476'

# Row 38
$ws.Range("F38").Style = "Normal"
$ws.Range("F38").Value = 'Doumassi,Kokou'
$ws.Range("H38").Value = 'Adjunct'
$ws.Range("I38").Value = '212-220-1476'

# Row 50
$ws.Range("F50").Style = "Normal"
$ws.Range("F50").Value = 'Anderson,Raheim Alan'
$ws.Range("H50").Value = 'Adjunct'
$ws.Range("I50").Value = '212-220-1476'

# Row 54
$ws.Range("F54").Style = "Normal"
$ws.Range("F54").Value = 'Santos,Jose Ramon R'
$ws.Range("H54").Value = 'Adjunct'
$ws.Range("I54").Value = '212-220-1476'

# Update active selection to reflect last edited cell
$ws.Range("I54").Select()
